# Updated cryptos list on Mon Jun 17 15:43:08 UTC 2024 with GitHub Actions
#
# Refresh of the Price / Volume(1h) columns scraped for this run, plus a
# rank swap between EnergySwap and Cosmos (rows 48-49 trade places).
# Every assignment is prefixed with a leading apostrophe so Excel stores
# numeric-looking text (e.g. "599.34", "1.00") as literal text instead of
# coercing it to a number -- matching the sheet's inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.405.38"
$ws.Range("E2").Value = "'  -1.84%  "
$ws.Range("D3").Value = "'3.503.20"
$ws.Range("E3").Value = "'  -2.32%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'599.34"
$ws.Range("E5").Value = "'  -1.63%  "
$ws.Range("D6").Value = "'142.25"
$ws.Range("E6").Value = "'  -3.56%  "
$ws.Range("D7").Value = "'3.502.35"
$ws.Range("E7").Value = "'  -2.31%  "
$ws.Range("E8").Value = "'  -0.27%  "
$ws.Range("E9").Value = "'  +5.54%  "
$ws.Range("E10").Value = "'  -3.09%  "
$ws.Range("D11").Value = "'7.81"
$ws.Range("E11").Value = "'  -2.47%  "
$ws.Range("E12").Value = "'  -3.14%  "
$ws.Range("D13").Value = "'4.105.80"
$ws.Range("E13").Value = "'  -2.13%  "
$ws.Range("D15").Value = "'28.24"
$ws.Range("E15").Value = "'  -6.11%  "
$ws.Range("D16").Value = "'3.510.52"
$ws.Range("E17").Value = "'  +1.33%  "
$ws.Range("D18").Value = "'65.394.74"
$ws.Range("D19").Value = "'10.81"
$ws.Range("E19").Value = "'  -5.46%  "
$ws.Range("E20").Value = "'  -2.50%  "
$ws.Range("E21").Value = "'  -5.04%  "
$ws.Range("D22").Value = "'417.35"
$ws.Range("E22").Value = "'  -3.45%  "
$ws.Range("E23").Value = "'  -5.17%  "
$ws.Range("D24").Value = "'76.91"
$ws.Range("E24").Value = "'  -2.76%  "
$ws.Range("D25").Value = "'3.649.17"
$ws.Range("E25").Value = "'  -2.20%  "
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E27").Value = "'  -5.94%  "
$ws.Range("E28").Value = "'  -3.13%  "
$ws.Range("D29").Value = "'8.88"
$ws.Range("E29").Value = "'  -4.53%  "
$ws.Range("E30").Value = "'  -5.97%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  +0.09%  "
$ws.Range("D32").Value = "'3.513.62"
$ws.Range("E32").Value = "'  -1.90%  "
$ws.Range("E33").Value = "'  -1.20%  "
$ws.Range("D34").Value = "'24.12"
$ws.Range("E34").Value = "'  -5.56%  "
$ws.Range("E35").Value = "'  -0.04%  "
$ws.Range("E36").Value = "'  -9.22%  "
$ws.Range("D37").Value = "'7.49"
$ws.Range("E37").Value = "'  -4.77%  "
$ws.Range("D38").Value = "'174.00"
$ws.Range("E38").Value = "'  +0.18%  "
$ws.Range("E39").Value = "'  -7.59%  "
$ws.Range("E40").Value = "'  -8.89%  "
$ws.Range("D41").Value = "'0.0807"
$ws.Range("E41").Value = "'  -5.78%  "
$ws.Range("E42").Value = "'  -4.93%  "
$ws.Range("D43").Value = "'0.852"
$ws.Range("E43").Value = "'  -4.98%  "
$ws.Range("D45").Value = "'1.75"
$ws.Range("E45").Value = "'  -8.38%  "
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "'  +0.13%  "
$ws.Range("E47").Value = "'  -8.39%  "
$ws.Range("B48").Value = "'Cosmos"
$ws.Range("C48").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "'7.01"
$ws.Range("E48").Value = "'  -2.80%  "
$ws.Range("B49").Value = "'EnergySwap"
$ws.Range("C49").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'23.05"
$ws.Range("E49").Value = "'  -3.78%  "
$ws.Range("E50").Value = "'  -8.57%  "
$ws.Range("D51").Value = "'0.899"
$ws.Range("E51").Value = "'  -5.41%  "
